# Update CVD (turnover / fill-rate) figures to reflect the refreshed source data files.
# Values below come from the regenerated BUR Testing_IPS.xlsx CVD export.
$wb = $excel.ActiveWorkbook

# --- IPS Clutches & Brakes Division ---
$ws = $wb.Worksheets.Item("IPS Clutches & Brakes Division")
$ws.Range("D2").Value = 0.0533
$ws.Range("D3").Value = 0.0533
$ws.Range("D4").Value = 0.0533
$ws.Range("I4").Value = 0.0197
$ws.Range("K4").Value = 0.0069
$ws.Range("L4").Value = 0.00691666666666667
$ws.Range("M4").Value = 0.02075
$ws.Range("N4").Value = 0.00691666666666667
$ws.Range("O4").Value = 0.00691666666666667
$ws.Range("P4").Value = 0.00691666666666667
$ws.Range("Q4").Value = 0.02075
$ws.Range("R4").Value = 0.00691666666666667
$ws.Range("S4").Value = 0.00691666666666667
$ws.Range("T4").Value = 0.00691666666666667
$ws.Range("U4").Value = 0.02075
$ws.Range("V4").Value = 0.083
$ws.Range("D5").Value = 0.459459459459459
$ws.Range("D6").Value = 0.459459459459459
$ws.Range("D7").Value = 0.459459459459459
$ws.Range("K7").Value = 0.5385
$ws.Range("L7").Value = 0.459459459459459
$ws.Range("M7").Value = 0.459459459459459
$ws.Range("N7").Value = 0.459459459459459
$ws.Range("O7").Value = 0.459459459459459
$ws.Range("P7").Value = 0.459459459459459
$ws.Range("Q7").Value = 0.459459459459459
$ws.Range("R7").Value = 0.459459459459459
$ws.Range("S7").Value = 0.459459459459459
$ws.Range("T7").Value = 0.459459459459459
$ws.Range("U7").Value = 0.459459459459459
$ws.Range("V7").Value = 0.459459459459459

# --- IPS Couplings Division ---
$ws = $wb.Worksheets.Item("IPS Couplings Division")
$ws.Range("D2").Value = 0.0533
$ws.Range("D3").Value = 0.0533
$ws.Range("D4").Value = 0.0533
$ws.Range("G4").Value = 0.0036
$ws.Range("I4").Value = 0.0263
$ws.Range("J4").Value = 0.006
$ws.Range("K4").Value = 0.0073
$ws.Range("L4").Value = 0.00794166666666667
$ws.Range("M4").Value = 0.023825
$ws.Range("N4").Value = 0.00794166666666667
$ws.Range("O4").Value = 0.00794166666666667
$ws.Range("P4").Value = 0.00794166666666667
$ws.Range("Q4").Value = 0.023825
$ws.Range("R4").Value = 0.00794166666666667
$ws.Range("S4").Value = 0.00794166666666667
$ws.Range("T4").Value = 0.00794166666666667
$ws.Range("U4").Value = 0.023825
$ws.Range("V4").Value = 0.0953
$ws.Range("D5").Value = 0.318181818181818
$ws.Range("D6").Value = 0.318181818181818
$ws.Range("D7").Value = 0.318181818181818
$ws.Range("J7").Value = 0.3333
$ws.Range("K7").Value = 0.3333
$ws.Range("L7").Value = 0.318181818181818
$ws.Range("M7").Value = 0.318181818181818
$ws.Range("N7").Value = 0.318181818181818
$ws.Range("O7").Value = 0.318181818181818
$ws.Range("P7").Value = 0.318181818181818
$ws.Range("Q7").Value = 0.318181818181818
$ws.Range("R7").Value = 0.318181818181818
$ws.Range("S7").Value = 0.318181818181818
$ws.Range("T7").Value = 0.318181818181818
$ws.Range("U7").Value = 0.318181818181818
$ws.Range("V7").Value = 0.318181818181818

# --- IPS Gearing Division ---
$ws = $wb.Worksheets.Item("IPS Gearing Division")
$ws.Range("D2").Value = 0.0533
$ws.Range("D3").Value = 0.0533
$ws.Range("D4").Value = 0.0533
$ws.Range("F4").Value = 0.0096
$ws.Range("H4").Value = 0.0113
$ws.Range("I4").Value = 0.0209
$ws.Range("J4").Value = 0.0113
$ws.Range("K4").Value = 0.0147
$ws.Range("L4").Value = 0.00935833333333333
$ws.Range("M4").Value = 0.028075
$ws.Range("N4").Value = 0.00935833333333333
$ws.Range("O4").Value = 0.00935833333333333
$ws.Range("P4").Value = 0.00935833333333333
$ws.Range("Q4").Value = 0.028075
$ws.Range("R4").Value = 0.00935833333333333
$ws.Range("S4").Value = 0.00935833333333333
$ws.Range("T4").Value = 0.00935833333333333
$ws.Range("U4").Value = 0.028075
$ws.Range("V4").Value = 0.1123
$ws.Range("D5").Value = 0.5
$ws.Range("D6").Value = 0.5
$ws.Range("D7").Value = 0.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.5
$ws.Range("N7").Value = 0.5
$ws.Range("O7").Value = 0.5
$ws.Range("P7").Value = 0.5
$ws.Range("Q7").Value = 0.5
$ws.Range("R7").Value = 0.5
$ws.Range("S7").Value = 0.5
$ws.Range("T7").Value = 0.5
$ws.Range("U7").Value = 0.5
$ws.Range("V7").Value = 0.5

# --- IPS Industrial Components Divi ---
$ws = $wb.Worksheets.Item("IPS Industrial Components Divi")
$ws.Range("D2").Value = 0.0533
$ws.Range("D3").Value = 0.0533
$ws.Range("D4").Value = 0.0533
$ws.Range("J4").Value = 0.0043
$ws.Range("K4").Value = 0.0197
$ws.Range("L4").Value = 0.00865833333333333
$ws.Range("M4").Value = 0.025975
$ws.Range("N4").Value = 0.00865833333333333
$ws.Range("O4").Value = 0.00865833333333333
$ws.Range("P4").Value = 0.00865833333333333
$ws.Range("Q4").Value = 0.025975
$ws.Range("R4").Value = 0.00865833333333333
$ws.Range("S4").Value = 0.00865833333333333
$ws.Range("T4").Value = 0.00865833333333333
$ws.Range("U4").Value = 0.025975
$ws.Range("V4").Value = 0.1039
$ws.Range("D5").Value = 0.526315789473684
$ws.Range("D6").Value = 0.526315789473684
$ws.Range("D7").Value = 0.526315789473684
$ws.Range("K7").Value = 0.5
$ws.Range("L7").Value = 0.526315789473684
$ws.Range("M7").Value = 0.526315789473684
$ws.Range("N7").Value = 0.526315789473684
$ws.Range("O7").Value = 0.526315789473684
$ws.Range("P7").Value = 0.526315789473684
$ws.Range("Q7").Value = 0.526315789473684
$ws.Range("R7").Value = 0.526315789473684
$ws.Range("S7").Value = 0.526315789473684
$ws.Range("T7").Value = 0.526315789473684
$ws.Range("U7").Value = 0.526315789473684
$ws.Range("V7").Value = 0.526315789473684

# --- IPS Segment Functions ---
$ws = $wb.Worksheets.Item("IPS Segment Functions")
$ws.Range("D2").Value = 0.0533
$ws.Range("D3").Value = 0.0533
$ws.Range("D4").Value = 0.0533
$ws.Range("F4").Value = 0.0084
$ws.Range("G4").Value = 0.0041
$ws.Range("H4").Value = 0.0109
$ws.Range("I4").Value = 0.0234
$ws.Range("J4").Value = 0.0067
$ws.Range("K4").Value = 0.0067
$ws.Range("L4").Value = 0.00735833333333333
$ws.Range("M4").Value = 0.022075
$ws.Range("N4").Value = 0.00735833333333333
$ws.Range("O4").Value = 0.00735833333333333
$ws.Range("P4").Value = 0.00735833333333333
$ws.Range("Q4").Value = 0.022075
$ws.Range("R4").Value = 0.00735833333333333
$ws.Range("S4").Value = 0.00735833333333333
$ws.Range("T4").Value = 0.00735833333333333
$ws.Range("U4").Value = 0.022075
$ws.Range("V4").Value = 0.0883
$ws.Range("D5").Value = 0.642857142857143
$ws.Range("D6").Value = 0.642857142857143
$ws.Range("D7").Value = 0.642857142857143
$ws.Range("G7").Value = 0.7273
$ws.Range("I7").Value = 0.697
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0.642857142857143
$ws.Range("M7").Value = 0.642857142857143
$ws.Range("N7").Value = 0.642857142857143
$ws.Range("O7").Value = 0.642857142857143
$ws.Range("P7").Value = 0.642857142857143
$ws.Range("Q7").Value = 0.642857142857143
$ws.Range("R7").Value = 0.642857142857143
$ws.Range("S7").Value = 0.642857142857143
$ws.Range("T7").Value = 0.642857142857143
$ws.Range("U7").Value = 0.642857142857143
$ws.Range("V7").Value = 0.642857142857143

# --- Integration ---
$ws = $wb.Worksheets.Item("Integration")
$ws.Range("D2").Value = 0.0533
$ws.Range("D3").Value = 0.0533

# --- L1_IPS ---
$ws = $wb.Worksheets.Item("L1_IPS")
$ws.Range("D2").Value = 0.0454
$ws.Range("D3").Value = 0.0454
$ws.Range("D4").Value = 0.0454
$ws.Range("J4").Value = 0.0073
$ws.Range("K4").Value = 0.0095
$ws.Range("L4").Value = 0.00781666666666667
$ws.Range("M4").Value = 0.02345
$ws.Range("N4").Value = 0.00781666666666667
$ws.Range("O4").Value = 0.00781666666666667
$ws.Range("P4").Value = 0.00781666666666667
$ws.Range("Q4").Value = 0.02345
$ws.Range("R4").Value = 0.00781666666666667
$ws.Range("S4").Value = 0.00781666666666667
$ws.Range("T4").Value = 0.00781666666666667
$ws.Range("U4").Value = 0.02345
$ws.Range("V4").Value = 0.0938
$ws.Range("D5").Value = 0.507246376811594
$ws.Range("D6").Value = 0.507246376811594
$ws.Range("D7").Value = 0.507246376811594
$ws.Range("G7").Value = 0.5526
$ws.Range("I7").Value = 0.5474
$ws.Range("J7").Value = 0.3913
$ws.Range("K7").Value = 0.45
$ws.Range("L7").Value = 0.507246376811594
$ws.Range("M7").Value = 0.507246376811594
$ws.Range("N7").Value = 0.507246376811594
$ws.Range("O7").Value = 0.507246376811594
$ws.Range("P7").Value = 0.507246376811594
$ws.Range("Q7").Value = 0.507246376811594
$ws.Range("R7").Value = 0.507246376811594
$ws.Range("S7").Value = 0.507246376811594
$ws.Range("T7").Value = 0.507246376811594
$ws.Range("U7").Value = 0.507246376811594
$ws.Range("V7").Value = 0.507246376811594
